$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (header banner text) ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Crime-complaint table value updates (rows 14-46) ---
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -50
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = -40
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -89.655172413793
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 19
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 72.727272727272
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 55.555555555555
$ws.Range("L15").Value = 133.333333333333
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = -26.315789473684
$ws.Range("C16").Value = 39
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = 8.333333333333
$ws.Range("F16").Value = 124
$ws.Range("G16").Value = 170
$ws.Range("H16").Value = -27.058823529411
$ws.Range("I16").Value = 90
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = -21.739130434782
$ws.Range("L16").Value = -10
$ws.Range("M16").Value = -27.419354838709
$ws.Range("N16").Value = -82.824427480916
$ws.Range("C17").Value = 42
$ws.Range("D17").Value = 48
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 207
$ws.Range("G17").Value = 215
$ws.Range("H17").Value = -3.720930232558
$ws.Range("I17").Value = 137
$ws.Range("J17").Value = 144
$ws.Range("K17").Value = -4.861111111111
$ws.Range("L17").Value = -2.142857142857
$ws.Range("M17").Value = 65.060240963855
$ws.Range("N17").Value = -48.106060606060
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = 154.545454545455
$ws.Range("F18").Value = 84
$ws.Range("G18").Value = 86
$ws.Range("H18").Value = -2.325581395348
$ws.Range("I18").Value = 58
$ws.Range("J18").Value = 56
$ws.Range("K18").Value = 3.571428571428
$ws.Range("L18").Value = -38.297872340425
$ws.Range("M18").Value = -36.956521739130
$ws.Range("N18").Value = -89.698046181172
$ws.Range("C19").Value = 96
$ws.Range("D19").Value = 105
$ws.Range("E19").Value = -8.571428571428
$ws.Range("F19").Value = 389
$ws.Range("G19").Value = 435
$ws.Range("H19").Value = -10.574712643678
$ws.Range("I19").Value = 264
$ws.Range("J19").Value = 286
$ws.Range("K19").Value = -7.692307692307
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 10.924369747899
$ws.Range("N19").Value = -54.794520547945
$ws.Range("C20").Value = 16
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 77.777777777777
$ws.Range("F20").Value = 50
$ws.Range("G20").Value = 76
$ws.Range("H20").Value = -34.210526315789
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = -17.073170731707
$ws.Range("L20").Value = -42.372881355932
$ws.Range("M20").Value = 47.826086956521
$ws.Range("N20").Value = -93.726937269372
$ws.Range("C21").Value = 224
$ws.Range("D21").Value = 214
$ws.Range("E21").Value = 4.672897196261
$ws.Range("F21").Value = 876
$ws.Range("G21").Value = 999
$ws.Range("H21").Value = -12.312312312312
$ws.Range("I21").Value = 600
$ws.Range("J21").Value = 656
$ws.Range("K21").Value = -8.536585365853
$ws.Range("L21").Value = -9.909909909909
$ws.Range("M21").Value = 5.263157894736
$ws.Range("N21").Value = -76.237623762376
$ws.Range("C22").Value = 3
$ws.Range("E22").Value = -62.5
$ws.Range("G22").Value = 31
$ws.Range("H22").Value = -58.064516129032
$ws.Range("I22").Value = 11
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = -50
$ws.Range("L22").Value = -15.384615384615
$ws.Range("M22").Value = -21.428571428571
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 16
$ws.Range("E23").Value = 93.75
$ws.Range("F23").Value = 93
$ws.Range("G23").Value = 99
$ws.Range("H23").Value = -6.060606060606
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = -11.428571428571
$ws.Range("L23").Value = 8.771929824561
$ws.Range("M23").Value = 67.567567567567
$ws.Range("C24").Value = 269
$ws.Range("D24").Value = 249
$ws.Range("E24").Value = 8.032128514056
$ws.Range("F24").Value = 969
$ws.Range("G24").Value = 854
$ws.Range("H24").Value = 13.466042154566
$ws.Range("I24").Value = 661
$ws.Range("J24").Value = 584
$ws.Range("K24").Value = 13.184931506849
$ws.Range("L24").Value = -7.032348804500
$ws.Range("M24").Value = 47.874720357941
$ws.Range("C25").Value = 137
$ws.Range("D25").Value = 146
$ws.Range("E25").Value = -6.164383561643
$ws.Range("F25").Value = 510
$ws.Range("G25").Value = 455
$ws.Range("H25").Value = 12.087912087912
$ws.Range("I25").Value = 346
$ws.Range("J25").Value = 321
$ws.Range("K25").Value = 7.788161993769
$ws.Range("L25").Value = -20.642201834862
$ws.Range("C26").Value = 86
$ws.Range("D26").Value = 63
$ws.Range("E26").Value = 36.507936507936
$ws.Range("F26").Value = 342
$ws.Range("G26").Value = 309
$ws.Range("H26").Value = 10.679611650485
$ws.Range("I26").Value = 240
$ws.Range("J26").Value = 205
$ws.Range("K26").Value = 17.073170731707
$ws.Range("L26").Value = 23.076923076923
$ws.Range("M26").Value = -7.335907335907
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = -20
$ws.Range("F27").Value = 22
$ws.Range("G27").Value = 17
$ws.Range("H27").Value = 29.411764705882
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = 30.769230769230
$ws.Range("L27").Value = 13.333333333333
$ws.Range("C28").Value = 12
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = 9.090909090909
$ws.Range("G28").Value = 42
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 24
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = -4
$ws.Range("L28").Value = -7.692307692307
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = -16.666666666666
$ws.Range("L29").Value = -44.444444444444
$ws.Range("M29").Value = -16.666666666666
$ws.Range("N29").Value = -91.071428571428
$ws.Range("C30").Value = 2
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -20
$ws.Range("L30").Value = -55.555555555555
$ws.Range("M30").Value = -20
$ws.Range("N30").Value = -92.156862745098
$ws.Range("D31").Value = 2
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = -90.909090909090
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = -80
$ws.Range("L31").Value = -75
$ws.Range("C33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 1
$ws.Range("L33").Value = -66.666666666666
$ws.Range("J39").Value = 48
$ws.Range("K39").Value = -38.461538461538
$ws.Range("L39").Value = -38.461538461538
$ws.Range("M39").Value = -84.313725490196
$ws.Range("N39").Value = -87.335092348285
$ws.Range("J42").Value = 3256
$ws.Range("K42").Value = 11.851597389213
$ws.Range("L42").Value = -14.853556485355
$ws.Range("M42").Value = -43.823326432022
$ws.Range("N42").Value = -46.552856204858
$ws.Range("J43").Value = 1402
$ws.Range("K43").Value = -50.248403122782
$ws.Range("L43").Value = -71.264603402336
$ws.Range("M43").Value = -87.745826413775
$ws.Range("N43").Value = -90.173126796102
$ws.Range("J44").Value = 6426
$ws.Range("K44").Value = 18.342541436464
$ws.Range("L44").Value = 10.336538461538
$ws.Range("M44").Value = -41.872455902306
$ws.Range("N44").Value = -53.431408073048
$ws.Range("J45").Value = 949
$ws.Range("K45").Value = -52.478718077115
$ws.Range("L45").Value = -69.805917912822
$ws.Range("M45").Value = -89.923550647695
$ws.Range("N45").Value = -92.110076488194
$ws.Range("J46").Value = 14252
$ws.Range("K46").Value = -15.181812771528
$ws.Range("L46").Value = -37.121680049413
$ws.Range("M46").Value = -70.465848806366
$ws.Range("N46").Value = -75.880043325209

# --- Cells that switched from blank/"N/A" text to a real number need the numeric style ---
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("I31").NumberFormat = "#,##0"
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("F33").NumberFormat = "#,##0"
$ws.Range("I33").NumberFormat = "#,##0"
